$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the test site credentials in row 3 (username/password pair)
$ws.Range("A3").Value = "mngr266814"
$ws.Range("B3").Value = "vAtarEt"
